$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale "AutoCorrelation Table" placeholder row (row 12) - this
# shifts rows 13:74 up by one.
$ws.Rows.Item(12).Delete()

# Append the completed autocorrelation-table test at the new last row (74).
$ws.Range("A74").Value = "Table_AutoCorrelation"
$ws.Range("B74").Value = "Test autocorrelation table"
$ws.Range("C74").Value = "table_autocorrelation_test"

# Restore the view state recorded after the edit.
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("F71").Select()
